$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.979.88"
$ws.Range("D3").Value = "1.845.22"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "2.112.97"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").Value = "1.842.86"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "35.011.95"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.50"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.123"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.29%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +25.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.96"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.749"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.26%  "
$ws.Range("E37").Value = "  +6.72%  "
$ws.Range("E38").Value = "  +11.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.01"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "1.348.14"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").Value = "2.032.53"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.43"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +18.13%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0669"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.35%  "
